# Refitting NCDEs to individual patients (for manuscript figure)
#
# Adds a new "Label" column (H) holding the true class label (0/Control,
# 1/MDD) used for each per-patient refit row, and refreshes the refitted
# loss / prediction / error values (columns D, E, F) that changed after
# re-running the per-patient NCDE fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell: H1 = "Label" ------------------------------------
# Match the bold / centered / thin-bordered style already used by the
# other header cells (B1:G1) and the row-label column (A2:A21) by copying
# the existing "Success %" header's format onto the new cell.
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats

# --- Refreshed metrics for the "100 iterations" patient block ---------
$ws.Range("D2").Value = [double]"5.997107583082974E-06"
$ws.Range("E2").Value = [double]"5.997107583082974E-06"
$ws.Range("H2").Value = 0

$ws.Range("D3").Value = 0.0002487910997451888
$ws.Range("E3").Value = 0.0002487910997451888
$ws.Range("H3").Value = 0

$ws.Range("D4").Value = 0.1772181380835149
$ws.Range("E4").Value = 0.1772181380835149
$ws.Range("H4").Value = 0

$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.4813067556821801
$ws.Range("E6").Value = 0.4813067556821801
$ws.Range("H6").Value = 0

$ws.Range("H7").Value = 1

$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 0.2478767395574011
$ws.Range("E9").Value = 0.7521232604425989
$ws.Range("H9").Value = 1

$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 0.4413685250936391
$ws.Range("E11").Value = 0.558631474906361
$ws.Range("F11").Value = 0.5372034311294556
$ws.Range("H11").Value = 1

# --- "Label" column for the "200 iterations" patient block ------------
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
